# "IDs de clases.xlsx" update (30/05/2019)
# The "Musica" (song/music) entries are renamed to "Cancion" (song) themed
# labels: Genero_Musica -> Genero_Cancion, Premios_Musica -> Premios_Cancion.
# Row order / numbering (column B) stays the same; only the text labels in
# column C for rows 15 and 16 change. The active selection is left on C16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C16 first, then C15, so that the new shared-string entries end up
# appended in the same order as the target workbook (Premios_Cancion before
# Genero_Cancion).
$ws.Range("C16").Value = "Premios_Cancion"
$ws.Range("C15").Value = "Genero_Cancion"

# Move the active selection to C16 (matching the saved cursor position).
$ws.Range("C16").Select()
